$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update #ALUTs values (column G) for rows 2 and 3
$ws.Range("G2").Value = 1612
$ws.Range("G3").Value = 1612

# Move the active selection to I5, matching the author's final cursor position
$ws.Range("I5").Select()
